$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.527.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.998.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.08%  "
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5023"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4237"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08955"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.116"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.11%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.997"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.43%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.977.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.466"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.014"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001114"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06811"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.937"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.548.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.312"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.316"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.311"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.060"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09948"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.557"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.835"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.805"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02459"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.260"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.99%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06400"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.32%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.297"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6557"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2044"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.01%  "
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6343"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.209"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.308"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.504"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000334"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06964"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.133"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.31%  "
